# Apply the edits described by the diff:
# 1. On the "Tables" sheet, rows 516-902, column H (SOURCE_FILE) changes from
#    "sample_padrao2_ello_nu_nordeste.xlsx" to
#    "sample_padrao2_ello_nu_nordeste_double_qtd - Copia.xlsm"
# 2. On the "Metadata" sheet:
#    - B4/B5: "3893" -> "4003"
#    - C4/C5: "LONDRINA AV-TIRADENTES" -> "ITAU RECIFE BOA VIAGEM"
#    - G5: "sample_padrao2_ello_nu_nordeste.xlsx" -> "sample_padrao2_ello_nu_nordeste_double_qtd - Copia.xlsm"

$wb = $excel.ActiveWorkbook

$oldSourceFile = "sample_padrao2_ello_nu_nordeste.xlsx"
$newSourceFile = "sample_padrao2_ello_nu_nordeste_double_qtd - Copia.xlsm"

# --- 1. Tables sheet: update SOURCE_FILE column (H) for rows 516-902 ---
$tablesWs = $wb.Worksheets.Item("Tables")

$startRow = 516
$endRow = 902
$col = 8  # column H

for ($r = $startRow; $r -le $endRow; $r++) {
    $cell = $tablesWs.Cells.Item($r, $col)
    if ($cell.Value2 -eq $oldSourceFile) {
        $cell.Value2 = $newSourceFile
    }
}

# --- 2. Metadata sheet: update agency number/name and source file ---
$metaWs = $wb.Worksheets.Item("Metadata")

# B4/B5 hold the agency number as text (e.g. "3893"); force text format so
# the numeric-looking string isn't silently coerced into a number.
$metaWs.Range("B4").NumberFormat = "@"
$metaWs.Range("B4").Value2 = "4003"
$metaWs.Range("C4").Value2 = "ITAU RECIFE BOA VIAGEM"

$metaWs.Range("B5").NumberFormat = "@"
$metaWs.Range("B5").Value2 = "4003"
$metaWs.Range("C5").Value2 = "ITAU RECIFE BOA VIAGEM"
$metaWs.Range("G5").Value2 = $newSourceFile
